$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A-T)
$rows = @(
  ,@("ECs","Myoc","Fzd7","ECs",2,0.6666666666666666,0.714093,2.142279,0.001553869939232348,0.001553869939232348,3,1,2.418393,7.255179,0.0919828589765645,0.0919828589765645,1.726957512549,15.542617612941,0.0001429293994883319,0.0001429293994883319)
  ,@("ECs","Myoc","Fzd7","FAPs",2,0.6666666666666666,0.714093,2.142279,0.001553869939232348,0.001553869939232348,3,1,10.11799233333333,30.353977,0.3848348311969811,0.3848348311969811,7.225187499287,65.026687493583,0.0005979832757665438,0.0005979832757665439)
  ,@("ECs","Myoc","Fzd7","sCs",2,0.6666666666666666,0.714093,2.142279,0.001553869939232348,0.001553869939232348,3,1,13.75539366666667,41.266181,0.5231823098264544,0.5231823098264544,9.822630329611,88.403672966499,0.000812957263977472,0.000812957263977472)
  ,@("FAPs","Myoc","Fzd7","ECs",3,1,453.4108886666666,1360.232666,0.9866243612803347,0.9866243612803348,3,1,2.418393,7.255179,0.0919828589765645,0.0919828589765645,1096.525719275246,9868.731473477214,0.09075252948649205,0.09075252948649205)
  ,@("FAPs","Myoc","Fzd7","FAPs",3,1,453.4108886666666,1360.232666,0.9866243612803347,0.9866243612803348,3,1,10.11799233333333,30.353977,0.3848348311969811,0.3848348311969811,4587.607895379187,41288.47105841268,0.3796874195281469,0.379687419528147)
  ,@("FAPs","Myoc","Fzd7","sCs",3,1,453.4108886666666,1360.232666,0.9866243612803347,0.9866243612803348,3,1,13.75539366666667,41.266181,0.5231823098264544,0.5231823098264544,6236.845266363172,56131.60739726855,0.5161844122656958,0.5161844122656958)
  ,@("sCs","Myoc","Fzd7","ECs",3,1,5.432785666666667,16.298357,0.01182176878043295,0.01182176878043295,3,1,2.418393,7.255179,0.0919828589765645,0.0919828589765645,13.138610826767,118.247497440903,0.001087400090584116,0.001087400090584116)
  ,@("sCs","Myoc","Fzd7","FAPs",3,1,5.432785666666667,16.298357,0.01182176878043295,0.01182176878043295,3,1,10.11799233333333,30.353977,0.3848348311969811,0.3848348311969811,54.96888372397656,494.719953515789,0.004549428393067653,0.004549428393067654)
  ,@("sCs","Myoc","Fzd7","sCs",3,1,5.432785666666667,16.298357,0.01182176878043295,0.01182176878043295,3,1,13.75539366666667,41.266181,0.5231823098264544,0.5231823098264544,74.73010555162412,672.5709499646171,0.006184940296781176,0.006184940296781176)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $rowData = $rows[$i]
  for ($c = 0; $c -lt $rowData.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
  }
}
